# Add a new "2021" data column (K) to the worksheet, mirroring the
# formatting already used in column J ("2020"), and update the
# worksheet's selection accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (styles/borders/number formats) from column J
# (rows 2-9) into the new column K so the new column matches the rest
# of the table's look.
$ws.Range("J2:J9").Copy()
$ws.Range("K2:K9").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the "2021" values for the new column.
$ws.Range("K3").Value = 2021
$ws.Range("K4").Value = 295
$ws.Range("K5").Value = 163
$ws.Range("K6").Value = 268
$ws.Range("K7").Value = 155
$ws.Range("K8").Value = 27
$ws.Range("K9").Value = 8

# Match the updated selection recorded in the workbook.
$ws.Range("L5").Select()
